# Update the "department" column (C) for every data row (2-13) on the
# "courses" sheet from "BRANSON SCHOOL OF BUSINESS AND TECHNOLOGY" to
# "Packages".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C13").Value = "Packages"

# Match the saved selection state from the source edit (active cell C13).
[void]$ws.Range("C13").Select()
